# Insert a new data row at row 210 (shifting existing rows 210..267 down to
# 211..268) and populate the new row with the Asterix / "1a (guarda)" entry
# for Región de La Araucanía.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 210:267 down by one to make room for the new record.
$ws.Rows("210:210").Insert()

# Populate the newly inserted row 210 with its values.
$ws.Cells.Item(210, 1).Value  = 11
$ws.Cells.Item(210, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(210, 3).Value  = "Bíobío"
$ws.Cells.Item(210, 4).Value  = 44736
$ws.Cells.Item(210, 5).Value  = 8
$ws.Cells.Item(210, 6).Value  = 100114001
$ws.Cells.Item(210, 7).Value  = "Papa"
$ws.Cells.Item(210, 8).Value  = "Asterix"
$ws.Cells.Item(210, 9).Value  = "1a (guarda)"
$ws.Cells.Item(210, 10).Value = 5000
$ws.Cells.Item(210, 11).Value = 8000
$ws.Cells.Item(210, 12).Value = 8500
$ws.Cells.Item(210, 13).Value = 8200
$ws.Cells.Item(210, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(210, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(210, 16).Value = 328
$ws.Cells.Item(210, 17).Value = 25
$ws.Cells.Item(210, 18).Value = "Hortaliza"
